# Product Backlog.xlsx update for Week 9 progress
# - "User stories" sheet: remove the "gameboard to be visible on the screen"
#   user story and renumber the remaining "No" column.
# - "Snake-game" sheet: add "Snake Logic" / "Gameboard Logic" / "Prey Logic"
#   as Level-1 task line items (abstraction work for snake/gameboard/prey),
#   restyle the "login form" row to match, and update selections.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("User stories")
$ws2 = $wb.Worksheets.Item("Snake-game")

# ---------------------------------------------------------------------------
# 1. "User stories" sheet — drop the row for
#    "the gameboard to be visible on the screen" (old row 13) and shift the
#    rest of the user stories up, renumbering the "No" column sequentially.
# ---------------------------------------------------------------------------
$ws1.Rows("13").Delete()

$ws1.Range("C13").Value = 5
$ws1.Range("C14").Value = 6
$ws1.Range("C16").Value = 7
$ws1.Range("C17").Value = 8
$ws1.Range("C18").Value = 9
$ws1.Range("C19").Value = 10
$ws1.Range("C20").Value = 11
$ws1.Range("C21").Value = 12
$ws1.Range("C23").Value = 13

# ---------------------------------------------------------------------------
# 2. "Snake-game" sheet — introduce abstraction-related Level-1 tasks:
#    Snake Logic, Gameboard Logic, Prey Logic (3 new rows under "login
#    form"), each weighted 0.1, matching the existing "login form" styling.
# ---------------------------------------------------------------------------
$ws2.Rows("37:39").Insert()

$ws2.Range("B36").Value = "login form"
$ws2.Range("F36").Value = 1

$ws2.Range("B37").Value = "Snake Logic"
$ws2.Range("F37").Value = 0.1

$ws2.Range("B38").Value = "Gameboard Logic"
$ws2.Range("F38").Value = 0.1

$ws2.Range("B39").Value = "Prey Logic"
$ws2.Range("F39").Value = 0.1

$taskRows = "B36", "B37", "B38", "B39"
foreach ($addr in $taskRows) {
    $cell = $ws2.Range($addr)
    $cell.Font.Bold = $true
    $cell.Font.Italic = $true
}

$weightRows = "F36", "F37", "F38", "F39"
foreach ($addr in $weightRows) {
    $cell = $ws2.Range($addr)
    $cell.NumberFormat = "0%"
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
}

# ---------------------------------------------------------------------------
# 3. Selections / active sheet — "User stories" becomes the active tab with
#    the cursor parked just below the last user story; "Snake-game" cursor
#    sits on the new last task row.
# ---------------------------------------------------------------------------
$ws2.Range("B39").Select()
$ws1.Activate()
$ws1.Range("D25").Select()
